$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "186,0,88"
$ws.Range("E3").Value = "186,0,88"
$ws.Range("E4").Value = "186,0,88"

$ws.Range("F7").Select()
